# update code tao report luong tai report co so
#
# Notion sync for "Lũy kế tháng SÓC TRĂNG" row 7 (Tháng 8) came back with a
# fresh last_edited_time / last_edited_by stamp plus the actual report
# numbers that used to be 0 / missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- last_edited_time ------------------------------------------------
# All rows whose last_edited_time currently reads 2024-07-31T18:24:00.000Z
# got bumped to the new Notion edit timestamp 2024-08-03T03:17:00.000Z -
# row 7 (which previously still had the older 2024-06-26 stamp) now joins
# that same "last edited" batch too.
$newEditedTime = "2024-08-03T03:17:00.000Z"
$ws.Range("D2").Value = $newEditedTime
$ws.Range("D3").Value = $newEditedTime
$ws.Range("D6").Value = $newEditedTime
$ws.Range("D7").Value = $newEditedTime
$ws.Range("D8").Value = $newEditedTime
$ws.Range("D11").Value = $newEditedTime
$ws.Range("D13").Value = $newEditedTime

# --- last_edited_by.id -------------------------------------------------
$ws.Range("N7").Value = "41cabcaf-915d-46a5-8eff-38727be27269"

# --- report numbers for Tháng 8 (row 7) --------------------------------
$ws.Range("S7").Value = 180000      # properties.Chi tiêu.number
$ws.Range("W7").Value = 5820000     # properties.Lũy kế.formula.number
$ws.Range("AA7").Value = 2000000    # properties.Dư nợ phát sinh.formula.number
$ws.Range("AE7").Value = 6000000    # properties.Tổng doanh thu.formula.number
$ws.Range("AH7").Value = 6000000    # properties.Đã thanh toán.number
$ws.Range("AK7").Value = 2          # properties.Số lượng đơn.number
$ws.Range("AN7").Value = 0          # properties.Thu nợ.number
$ws.Range("AQ7").Value = 8000000    # properties.Đơn giá.number
